# Swap the contents of columns C ("codeforiati:group-name") and D
# ("codeforiati:group-code") for every row in the used range, including
# the header row. This mirrors the shared-string reordering in the
# source diff, whose net effect on the worksheet is that the "name"
# and "code" columns trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Text
    $dVal = $dCell.Text

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
